$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Priority column (D) values for the user stories
$ws.Range("D4").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("D8").Value = 3
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("D12").Value = 3
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 3
$ws.Range("D15").Value = 3

# Move the active selection to D16
$ws.Range("D16").Select() | Out-Null
